$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on the Price and Volume columns so that values such as
# "65.517.12", "0.999", "7.05" etc. are preserved exactly as typed instead of
# being auto-converted into numeric cells by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '65.517.12'
$ws.Range("E2").Value = '  -0.36%  '

# Row 3
$ws.Range("D3").Value = '3.564.51'
$ws.Range("E3").Value = '  +3.07%  '

# Row 4
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.05%  '

# Row 5
$ws.Range("D5").Value = '602.26'
$ws.Range("E5").Value = '  +0.87%  '

# Row 6
$ws.Range("D6").Value = '140.83'
$ws.Range("E6").Value = '  +2.55%  '

# Row 7
$ws.Range("D7").Value = '3.563.77'
$ws.Range("E7").Value = '  +3.12%  '

# Row 8
$ws.Range("E8").Value = '  +0.07%  '

# Row 9
$ws.Range("D9").Value = '0.494'
$ws.Range("E9").Value = '  -0.10%  '

# Row 10
$ws.Range("E10").Value = '  +2.96%  '

# Row 11
$ws.Range("D11").Value = '7.05'
$ws.Range("E11").Value = '  -6.07%  '

# Row 12
$ws.Range("D12").Value = '0.397'
$ws.Range("E12").Value = '  +4.22%  '

# Row 13
$ws.Range("D13").Value = '4.170.83'
$ws.Range("E13").Value = '  +3.18%  '

# Row 14
$ws.Range("E14").Value = '  +2.61%  '

# Row 15
$ws.Range("D15").Value = '27.24'
$ws.Range("E15").Value = '  +2.06%  '

# Row 16
$ws.Range("D16").Value = '3.564.13'
$ws.Range("E16").Value = '  +4.17%  '

# Row 17
$ws.Range("E17").Value = '  +1.48%  '

# Row 18
$ws.Range("D18").Value = '65.569.81'
$ws.Range("E18").Value = '  -0.18%  '

# Row 19
$ws.Range("D19").Value = '10.38'
$ws.Range("E19").Value = '  +4.83%  '

# Row 20
$ws.Range("D20").Value = '5.92'
$ws.Range("E20").Value = '  +2.17%  '

# Row 21
$ws.Range("D21").Value = '14.29'
$ws.Range("E21").Value = '  +3.76%  '

# Row 22
$ws.Range("D22").Value = '397.94'
$ws.Range("E22").Value = '  +0.46%  '

# Row 23
$ws.Range("D23").Value = '0.575'
$ws.Range("E23").Value = '  +4.50%  '

# Row 24
$ws.Range("D24").Value = '3.705.76'
$ws.Range("E24").Value = '  +2.84%  '

# Row 25
$ws.Range("D25").Value = '74.43'
$ws.Range("E25").Value = '  +1.12%  '

# Row 26
$ws.Range("E26").Value = '  -0.01%  '

# Row 27
$ws.Range("D27").Value = '0.0000118'
$ws.Range("E27").Value = '  +9.55%  '

# Row 28
$ws.Range("D28").Value = '7.92'
$ws.Range("E28").Value = '  +9.30%  '

# Row 29
$ws.Range("E29").Value = '  -0.11%  '

# Row 30
$ws.Range("E30").Value = '  +0.56%  '

# Row 31
$ws.Range("D31").Value = '8.36'
$ws.Range("E31").Value = '  +0.75%  '

# Row 32
$ws.Range("D32").Value = '3.582.06'
$ws.Range("E32").Value = '  +3.45%  '

# Row 33
$ws.Range("B33").Value = 'EthereumClassic'
$ws.Range("C33").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D33").Value = '23.96'
$ws.Range("E33").Value = '  +4.09%  '

# Row 34
$ws.Range("B34").Value = 'Kaspa'
$ws.Range("C34").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D34").Value = '0.148'
$ws.Range("E34").Value = '  +0.80%  '

# Row 35
$ws.Range("B35").Value = 'USDe'
$ws.Range("C35").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  +0.02%  '

# Row 36
$ws.Range("D36").Value = '1.28'
$ws.Range("E36").Value = '  +4.33%  '

# Row 37
$ws.Range("D37").Value = '7.11'
$ws.Range("E37").Value = '  +2.33%  '

# Row 38
$ws.Range("B38").Value = 'ImmutableX'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D38").Value = '1.55'
$ws.Range("E38").Value = '  +1.37%  '

# Row 39
$ws.Range("B39").Value = 'Monero'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D39").Value = '168.19'
$ws.Range("E39").Value = '  -2.95%  '

# Row 40
$ws.Range("D40").Value = '5.06'
$ws.Range("E40").Value = '  +4.19%  '

# Row 41
$ws.Range("D41").Value = '0.0809'
$ws.Range("E41").Value = '  +3.71%  '

# Row 42
$ws.Range("D42").Value = '0.836'
$ws.Range("E42").Value = '  +1.46%  '

# Row 43
$ws.Range("D43").Value = '26.89'
$ws.Range("E43").Value = '  +14.99%  '

# Row 44
$ws.Range("D44").Value = '42.91'

# Row 45
$ws.Range("D45").Value = '0.999'
$ws.Range("E45").Value = '  -0.02%  '

# Row 46
$ws.Range("D46").Value = '4.46'
$ws.Range("E46").Value = '  -0.02%  '

# Row 47
$ws.Range("D47").Value = '1.70'
$ws.Range("E47").Value = '  +3.52%  '

# Row 48
$ws.Range("E48").Value = '  +7.87%  '

# Row 49
$ws.Range("D49").Value = '2.451.84'
$ws.Range("E49").Value = '  +10.41%  '

# Row 50
$ws.Range("D50").Value = '6.85'
$ws.Range("E50").Value = '  +3.88%  '

# Row 51
$ws.Range("B51").Value = 'dogwifhat'
$ws.Range("C51").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D51").Value = '2.14'
$ws.Range("E51").Value = '  +0.46%  '
